$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D), Volume(1h) (E), and Coin/Link (B/C) refresh scraped from coinranking.com
# D-column cells are forced to text format so numeric-looking strings (e.g. "385.20")
# keep their exact textual representation instead of being parsed as numbers.
$textUpdates = @(
    ,@("D2", "51.468.72")
    ,@("D3", "3.098.62")
    ,@("D5", "385.20")
    ,@("D6", "103.19")
    ,@("D9", "0.583")
    ,@("D10", "36.92")
    ,@("D12", "0.0853")
    ,@("D13", "3.592.75")
    ,@("D14", "18.51")
    ,@("D15", "7.80")
    ,@("D16", "3.101.86")
    ,@("D17", "0.989")
    ,@("D18", "10.90")
    ,@("D19", "51.522.69")
    ,@("D20", "3.24")
    ,@("D21", "12.37")
    ,@("D22", "0.0₃0961")
    ,@("D23", "69.75")
    ,@("D24", "266.44")
    ,@("D26", "8.06")
    ,@("D27", "26.92")
    ,@("D28", "1.00")
    ,@("D29", "7.14")
    ,@("D32", "10.33")
    ,@("D34", "35.10")
    ,@("D38", "3.34")
    ,@("D39", "0.290")
    ,@("D40", "1.87")
    ,@("D41", "128.90")
    ,@("D43", "16.52")
    ,@("D46", "22.11")
    ,@("D47", "2.54")
    ,@("D48", "2.08")
    ,@("D49", "2.064.03")
    ,@("D50", "0.939")
    ,@("D51", "0.0327")
)

foreach ($pair in $textUpdates) {
    $cell = $ws.Range($pair[0])
    $cell.NumberFormat = "@"
    $cell.Value = $pair[1]
    $cell.Style = "Normal"
}

# B, C, E columns are naturally non-numeric text; plain assignment is sufficient.
$plainUpdates = @(
    ,@("E2", "  -0.27%  ")
    ,@("E3", "  +2.50%  ")
    ,@("E4", "  +0.04%  ")
    ,@("E5", "  +1.68%  ")
    ,@("E6", "  +0.21%  ")
    ,@("E7", "  -1.32%  ")
    ,@("E8", "  +0.02%  ")
    ,@("E9", "  -1.73%  ")
    ,@("E10", "  +0.94%  ")
    ,@("E11", "  +0.07%  ")
    ,@("E12", "  -0.79%  ")
    ,@("E13", "  +2.56%  ")
    ,@("E14", "  +0.04%  ")
    ,@("E15", "  +0.87%  ")
    ,@("E16", "  +2.42%  ")
    ,@("E17", "  +1.24%  ")
    ,@("E18", "  +2.80%  ")
    ,@("E19", "  -0.20%  ")
    ,@("E20", "  +6.49%  ")
    ,@("E21", "  -0.64%  ")
    ,@("E22", "  -0.01%  ")
    ,@("E23", "  -0.23%  ")
    ,@("E24", "  -0.78%  ")
    ,@("E25", "  +0.20%  ")
    ,@("E26", "  -2.23%  ")
    ,@("E27", "  +2.34%  ")
    ,@("B28", "Dai")
    ,@("C28", "https://coinranking.com/coin/MoTuySvg7+dai-dai")
    ,@("E28", "  -0.09%  ")
    ,@("B29", "RenderToken")
    ,@("C29", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr")
    ,@("E29", "  -7.15%  ")
    ,@("E30", "  -2.97%  ")
    ,@("E31", "  -2.44%  ")
    ,@("E32", "  +0.52%  ")
    ,@("E33", "  +5.30%  ")
    ,@("E34", "  +2.83%  ")
    ,@("E35", "  +0.46%  ")
    ,@("E36", "  -1.30%  ")
    ,@("E37", "  -0.19%  ")
    ,@("E38", "  +1.03%  ")
    ,@("E39", "  +0.40%  ")
    ,@("E40", "  +0.54%  ")
    ,@("E41", "  +1.69%  ")
    ,@("E42", "  -0.01%  ")
    ,@("E43", "  -3.43%  ")
    ,@("E44", "  -3.52%  ")
    ,@("E45", "  -0.30%  ")
    ,@("E46", "  +1.72%  ")
    ,@("E47", "  +6.66%  ")
    ,@("E48", "  -0.28%  ")
    ,@("E49", "  +1.60%  ")
    ,@("E50", "  +19.32%  ")
    ,@("E51", "  +1.58%  ")
)

foreach ($pair in $plainUpdates) {
    $ws.Range($pair[0]).Value = $pair[1]
}
